# Update the workbook per the commit:
#   - bump URL from ibm.com -> linuxforhealth.org
#   - bump Version 7.0.0 -> 8.0.0
#   - update Date and Publisher
#   - fix a data-alignment issue on the "Elements" sheet where the
#     ele-1/ext-1 Constraint(s) text was attached to the "Extension" row
#     instead of the "Extension.extension" row, and the Fixed Value on
#     the "Extension.url" row still referenced the old URL

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metadata sheet
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-job-title"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---------------------------------------------------------------------
# Elements sheet
# ---------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# The ele-1 / ext-1 constraint text was incorrectly shown on the
# "Extension" row (row 2); it actually belongs on the
# "Extension.extension" row (row 4), which already carries it. Clear
# the erroneous copy from row 2.
$elements.Range("AI2").Value = ""

# The Fixed Value of Extension.url (row 5) should track the new URL.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-job-title"
